$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$whatIDid = "Watched children's shows from my childhood and their continuations and spin-offs, and read simple manga."

# Add a new week 5 row of data.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1.667939814814815
$ws.Range("B6").NumberFormat = "[h]:mm:ss"
$ws.Range("C6").Value = "Your Lie in April  (Text with visuals, Japanese, Familiar):32; Drake and Josh (Audiovisual, English, Familiar):30;"
$ws.Range("D6").Value = $whatIDid

# Correct existing text in D5: remove the stray "with" before the comma.
$ws.Range("D5").Value = $whatIDid

$ws.Range("D6").Select()
